# "Section conseils et astuces"
# Set D7 ("État de la doc" sheet) status to "Terminé" and update the
# active selection to D7 (removing the scrolled topLeftCell).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("État de la doc")

$ws.Range("D7").Value = "Terminé"

$ws.Activate()
$ws.Range("D7").Select()
